$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.072359800338745
$ws.Range("B1").Value = 2.417393684387207
$ws.Range("C1").Value = 5.0977783203125
$ws.Range("D1").Value = 2.294564008712769
$ws.Range("E1").Value = 1.302074432373047
